# Add a third "type" column classifying each translation key/value row as
# a "title" or "description" string (and the header row as "type").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "type"
$ws.Range("C2").Value = "title"
$ws.Range("C3").Value = "title"
$ws.Range("C4").Value = "description"
$ws.Range("C5").Value = "description"
$ws.Range("C6").Value = "description"

# Match the re-saved column widths for the (now wider) text columns.
$ws.Columns.Item(1).ColumnWidth = 28.6666666666667
$ws.Columns.Item(2).ColumnWidth = 119.6666666666667

# Leave the selection where the workbook was left after the edit.
$ws.Range("C7").Select()
